$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the data row for account 002687737 / JOSE / 4421.99 (Excel row 6)
$ws.Rows.Item(6).Delete()
